$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title row: new subtitle text, taller row
$ws.Range("A1").Value = "Arbeitsjournal Semesterarbeit, Modul 2"
$ws.Rows.Item(1).RowHeight = 21

# Row 6: updated date / topic / hours
$ws.Range("A6").Value = 43473
$ws.Range("B6").Value = "Besprechung"
$ws.Range("C6").Value = 1

# Row 7: updated date / topic / hours
$ws.Range("A7").Value = 43485
$ws.Range("B7").Value = "Stored Procedure BillPoD"
$ws.Range("C7").Value = 4

# Rows 8-19: clear out the old journal entries (keep formatting/styles)
$ws.Range("A8:C19").ClearContents()

# Selection / view state as left by the author
[void]$ws.Range("A8:C22").Select()
